# Auto-generated Word COM-interop edit script.
# Updates the worksheet date heading and the 100 addition/subtraction
# problem answers in the single table, cell by cell, in document order.

$d = $word.ActiveDocument

# --- Update the date/weekday heading paragraph -----------------------
$d.Content.Find.Execute("2026-01-21 Wednesday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2026-01-22 Thursday", 2) | Out-Null

# --- Update each table cell's answer text -----------------------------
$newValues = @(
    "56-42=14",
    "67-23=44",
    "93-0=93",
    "3+37=40",
    "22-17=5",
    "14+47=61",
    "81-29=52",
    "92-56=36",
    "11-10=1",
    "29+62=91",
    "98-22=76",
    "53+45=98",
    "22+25=47",
    "41+36=77",
    "94-60=34",
    "99-77=22",
    "98-33=65",
    "26+36=62",
    "64-13=51",
    "47+32=79",
    "70-25=45",
    "25+74=99",
    "28+67=95",
    "47-44=3",
    "50+17=67",
    "36-12=24",
    "46+19=65",
    "88-58=30",
    "34+5=39",
    "12+24=36",
    "12+26=38",
    "11+23=34",
    "35+4=39",
    "30+15=45",
    "18+7=25",
    "85-16=69",
    "35-33=2",
    "44-14=30",
    "89-45=44",
    "95-10=85",
    "23+16=39",
    "69-23=46",
    "19+62=81",
    "59-2=57",
    "66-34=32",
    "79-1=78",
    "74-61=13",
    "66-40=26",
    "67-13=54",
    "39-22=17",
    "77-19=58",
    "3+66=69",
    "91-23=68",
    "72-57=15",
    "21+72=93",
    "71-44=27",
    "5+22=27",
    "30+4=34",
    "3+2=5",
    "80-0=80",
    "38+5=43",
    "89+0=89",
    "91-70=21",
    "60+8=68",
    "76-0=76",
    "61-54=7",
    "45+45=90",
    "85-45=40",
    "22+60=82",
    "68-29=39",
    "49-40=9",
    "62+30=92",
    "96-61=35",
    "50+3=53",
    "83-32=51",
    "61-15=46",
    "63-2=61",
    "16+36=52",
    "81-8=73",
    "56-34=22",
    "69+0=69",
    "18+53=71",
    "14+31=45",
    "87-69=18",
    "90-20=70",
    "93-51=42",
    "3+14=17",
    "26+61=87",
    "66+14=80",
    "76-75=1",
    "94-91=3",
    "49+40=89",
    "90-68=22",
    "12+12=24",
    "81-53=28",
    "32+8=40",
    "91-35=56",
    "41+41=82",
    "95-60=35",
    "25+49=74"
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

if (($rows * $cols) -ne $newValues.Length) {
    Write-Host "WARNING: table has" ($rows * $cols) "cells but" $newValues.Length "replacement values were supplied."
}

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated" $idx "table cells."
